# === Auto-generated Excel COM-interop edit script ===
# Implements: add "Player Info" + "ODI Batting Extra" sheets; rename
# MATCH_CARD_LINK -> MATCH_CODE and replace URL values with bare match
# codes on "ODI Batting" / "ODI Bowling"; drop now-empty INNING_NUMBER cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Insert "Player Info" before the first sheet, "ODI Batting Extra"
#    after the last sheet, so the final order is:
#    Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ---------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

# ---------------------------------------------------------------
# Helper: write a value to a cell, forcing TEXT storage (so purely
# numeric-looking strings like match codes do not get silently
# coerced to Excel numbers), without leaving stray cell styles.
# ---------------------------------------------------------------
function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

function Set-NumberCell($sheet, $row, $col, $num) {
    $sheet.Cells.Item($row, $col).Value = $num
}

function Clear-Cell($sheet, $row, $col) {
    $sheet.Cells.Item($row, $col).ClearContents()
}

# ---------------------------------------------------------------
# 2. Populate "Player Info"
# ---------------------------------------------------------------
Set-TextCell $playerInfo 1 1 'ID'
Set-TextCell $playerInfo 1 2 'NAME'
Set-TextCell $playerInfo 1 3 'BATTING_HAND'
Set-TextCell $playerInfo 1 4 'BOWL_STYLE'
Set-TextCell $playerInfo 2 1 '4328'
Set-TextCell $playerInfo 2 2 'Craig Alexander Young'
Set-TextCell $playerInfo 2 3 'Right Handed'
Set-TextCell $playerInfo 2 4 'Right Arm Medium'

# ---------------------------------------------------------------
# 3. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE header + values
# ---------------------------------------------------------------
Set-TextCell $battingSheet 1 4 'MATCH_CODE'
Set-TextCell $battingSheet 2 4 '3674'
Set-TextCell $battingSheet 3 4 '3675'
Set-TextCell $battingSheet 4 4 '3676'
Set-TextCell $battingSheet 5 4 '3722'
Set-TextCell $battingSheet 6 4 '3724'
Set-TextCell $battingSheet 7 4 '3730'
Set-TextCell $battingSheet 8 4 '3800'
Set-TextCell $battingSheet 9 4 '3829'
Set-TextCell $battingSheet 10 4 '3842'
Set-TextCell $battingSheet 11 4 '3935'
Set-TextCell $battingSheet 12 4 '4000'
Set-TextCell $battingSheet 13 4 '4003'
Set-TextCell $battingSheet 14 4 '4026'
Set-TextCell $battingSheet 15 4 '4397'
Set-TextCell $battingSheet 16 4 '4426'
Set-TextCell $battingSheet 17 4 '4427'
Set-TextCell $battingSheet 18 4 '4428'
Set-TextCell $battingSheet 19 4 '4439'
Set-TextCell $battingSheet 20 4 '4442'
Set-TextCell $battingSheet 21 4 '4444'
Set-TextCell $battingSheet 22 4 '4446'
Set-TextCell $battingSheet 23 4 '4448'
Set-TextCell $battingSheet 24 4 '4466'
Set-TextCell $battingSheet 25 4 '4467'
Set-TextCell $battingSheet 26 4 '4468'
Set-TextCell $battingSheet 27 4 '4474'
Set-TextCell $battingSheet 28 4 '4475'
Set-TextCell $battingSheet 29 4 '4478'
Set-TextCell $battingSheet 30 4 '4492'
Set-TextCell $battingSheet 31 4 '4494'
Set-TextCell $battingSheet 32 4 '4519'
Set-TextCell $battingSheet 33 4 '4520'
Set-TextCell $battingSheet 34 4 '4522'
Set-TextCell $battingSheet 35 4 '4605'
Set-TextCell $battingSheet 36 4 '4608'
Set-TextCell $battingSheet 37 4 '4614'

# "ODI Batting": drop the now-empty INNING_NUMBER cells (col B)
Clear-Cell $battingSheet 2 2
Clear-Cell $battingSheet 3 2
Clear-Cell $battingSheet 5 2
Clear-Cell $battingSheet 6 2
Clear-Cell $battingSheet 8 2
Clear-Cell $battingSheet 10 2
Clear-Cell $battingSheet 13 2
Clear-Cell $battingSheet 18 2
Clear-Cell $battingSheet 19 2
Clear-Cell $battingSheet 20 2
Clear-Cell $battingSheet 25 2
Clear-Cell $battingSheet 27 2
Clear-Cell $battingSheet 28 2
Clear-Cell $battingSheet 33 2

# ---------------------------------------------------------------
# 4. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE header + values
# ---------------------------------------------------------------
Set-TextCell $bowlingSheet 1 2 'MATCH_CODE'
Set-TextCell $bowlingSheet 2 2 '3674'
Set-TextCell $bowlingSheet 3 2 '3675'
Set-TextCell $bowlingSheet 4 2 '3676'
Set-TextCell $bowlingSheet 5 2 '3722'
Set-TextCell $bowlingSheet 6 2 '3724'
Set-TextCell $bowlingSheet 7 2 '3730'
Set-TextCell $bowlingSheet 8 2 '3829'
Set-TextCell $bowlingSheet 9 2 '3842'
Set-TextCell $bowlingSheet 10 2 '3935'
Set-TextCell $bowlingSheet 11 2 '4000'
Set-TextCell $bowlingSheet 12 2 '4003'
Set-TextCell $bowlingSheet 13 2 '4026'
Set-TextCell $bowlingSheet 14 2 '4397'
Set-TextCell $bowlingSheet 15 2 '4426'
Set-TextCell $bowlingSheet 16 2 '4427'
Set-TextCell $bowlingSheet 17 2 '4428'
Set-TextCell $bowlingSheet 18 2 '4439'
Set-TextCell $bowlingSheet 19 2 '4442'
Set-TextCell $bowlingSheet 20 2 '4444'
Set-TextCell $bowlingSheet 21 2 '4446'
Set-TextCell $bowlingSheet 22 2 '4448'
Set-TextCell $bowlingSheet 23 2 '4466'
Set-TextCell $bowlingSheet 24 2 '4467'
Set-TextCell $bowlingSheet 25 2 '4468'
Set-TextCell $bowlingSheet 26 2 '4475'
Set-TextCell $bowlingSheet 27 2 '4478'
Set-TextCell $bowlingSheet 28 2 '4492'
Set-TextCell $bowlingSheet 29 2 '4519'
Set-TextCell $bowlingSheet 30 2 '4520'
Set-TextCell $bowlingSheet 31 2 '4522'
Set-TextCell $bowlingSheet 32 2 '4605'
Set-TextCell $bowlingSheet 33 2 '4608'
Set-TextCell $bowlingSheet 34 2 '4614'

# ---------------------------------------------------------------
# 5. Populate "ODI Batting Extra"
# ---------------------------------------------------------------
Set-TextCell $battingExtra 1 1 'MATCH_CODE'
Set-TextCell $battingExtra 1 2 'BATTING_POSITION'
Set-TextCell $battingExtra 1 3 'NUM_4'
Set-TextCell $battingExtra 1 4 'NUM_6'
Set-TextCell $battingExtra 1 5 'PERCENT_RUNS_OF_TOTAL'
Set-TextCell $battingExtra 1 6 'MAN_OF_MATCH'
Set-TextCell $battingExtra 2 1 '4428'
Set-NumberCell $battingExtra 2 2 10
Set-TextCell $battingExtra 2 3 ''
Set-TextCell $battingExtra 2 4 ''
Set-TextCell $battingExtra 2 5 ''
Set-TextCell $battingExtra 2 6 'NO'
Set-TextCell $battingExtra 3 1 '4439'
Set-TextCell $battingExtra 3 2 ''
Set-TextCell $battingExtra 3 3 ''
Set-TextCell $battingExtra 3 4 ''
Set-TextCell $battingExtra 3 5 ''
Set-TextCell $battingExtra 3 6 'NO'
Set-TextCell $battingExtra 4 1 '4442'
Set-NumberCell $battingExtra 4 2 11
Set-TextCell $battingExtra 4 3 ''
Set-TextCell $battingExtra 4 4 ''
Set-TextCell $battingExtra 4 5 ''
Set-TextCell $battingExtra 4 6 'NO'
Set-TextCell $battingExtra 5 1 '4444'
Set-NumberCell $battingExtra 5 2 11
Set-TextCell $battingExtra 5 3 '0'
Set-TextCell $battingExtra 5 4 '0'
Set-TextCell $battingExtra 5 5 ''
Set-TextCell $battingExtra 5 6 'NO'
Set-TextCell $battingExtra 6 1 '4446'
Set-NumberCell $battingExtra 6 2 10
Set-TextCell $battingExtra 6 3 '0'
Set-TextCell $battingExtra 6 4 '0'
Set-TextCell $battingExtra 6 5 ''
Set-TextCell $battingExtra 6 6 'NO'
Set-TextCell $battingExtra 7 1 '4448'
Set-NumberCell $battingExtra 7 2 10
Set-TextCell $battingExtra 7 3 ''
Set-TextCell $battingExtra 7 4 ''
Set-TextCell $battingExtra 7 5 ''
Set-TextCell $battingExtra 7 6 'NO'
Set-TextCell $battingExtra 8 1 '4466'
Set-NumberCell $battingExtra 8 2 11
Set-TextCell $battingExtra 8 3 '0'
Set-TextCell $battingExtra 8 4 '0'
Set-TextCell $battingExtra 8 5 '0.61%'
Set-TextCell $battingExtra 8 6 'NO'
Set-TextCell $battingExtra 9 1 '4467'
Set-TextCell $battingExtra 9 2 ''
Set-TextCell $battingExtra 9 3 ''
Set-TextCell $battingExtra 9 4 ''
Set-TextCell $battingExtra 9 5 ''
Set-TextCell $battingExtra 9 6 'NO'
Set-TextCell $battingExtra 10 1 '4468'
Set-NumberCell $battingExtra 10 2 11
Set-TextCell $battingExtra 10 3 ''
Set-TextCell $battingExtra 10 4 ''
Set-TextCell $battingExtra 10 5 ''
Set-TextCell $battingExtra 10 6 'NO'
Set-TextCell $battingExtra 11 1 '4474'
Set-TextCell $battingExtra 11 2 ''
Set-TextCell $battingExtra 11 3 ''
Set-TextCell $battingExtra 11 4 ''
Set-TextCell $battingExtra 11 5 ''
Set-TextCell $battingExtra 11 6 'NO'
Set-TextCell $battingExtra 12 1 '4475'
Set-NumberCell $battingExtra 12 2 10
Set-TextCell $battingExtra 12 3 '0'
Set-TextCell $battingExtra 12 4 '0'
Set-TextCell $battingExtra 12 5 '0.88%'
Set-TextCell $battingExtra 12 6 'NO'
Set-TextCell $battingExtra 13 1 '4478'
Set-TextCell $battingExtra 13 2 ''
Set-TextCell $battingExtra 13 3 ''
Set-TextCell $battingExtra 13 4 ''
Set-TextCell $battingExtra 13 5 ''
Set-TextCell $battingExtra 13 6 'NO'
Set-TextCell $battingExtra 14 1 '4492'
Set-TextCell $battingExtra 14 2 ''
Set-TextCell $battingExtra 14 3 ''
Set-TextCell $battingExtra 14 4 ''
Set-TextCell $battingExtra 14 5 ''
Set-TextCell $battingExtra 14 6 'NO'
Set-TextCell $battingExtra 15 1 '4494'
Set-NumberCell $battingExtra 15 2 10
Set-TextCell $battingExtra 15 3 ''
Set-TextCell $battingExtra 15 4 ''
Set-TextCell $battingExtra 15 5 ''
Set-TextCell $battingExtra 15 6 'NO'
Set-TextCell $battingExtra 16 1 '4519'
Set-TextCell $battingExtra 16 2 ''
Set-TextCell $battingExtra 16 3 ''
Set-TextCell $battingExtra 16 4 ''
Set-TextCell $battingExtra 16 5 ''
Set-TextCell $battingExtra 16 6 'NO'
Set-TextCell $battingExtra 17 1 '4520'
Set-NumberCell $battingExtra 17 2 10
Set-TextCell $battingExtra 17 3 '0'
Set-TextCell $battingExtra 17 4 '0'
Set-TextCell $battingExtra 17 5 ''
Set-TextCell $battingExtra 17 6 'NO'
Set-TextCell $battingExtra 18 1 '4522'
Set-NumberCell $battingExtra 18 2 10
Set-TextCell $battingExtra 18 3 '0'
Set-TextCell $battingExtra 18 4 '0'
Set-TextCell $battingExtra 18 5 ''
Set-TextCell $battingExtra 18 6 'NO'
Set-TextCell $battingExtra 19 1 '4605'
Set-NumberCell $battingExtra 19 2 10
Set-TextCell $battingExtra 19 3 '1'
Set-TextCell $battingExtra 19 4 '0'
Set-TextCell $battingExtra 19 5 '1.67%'
Set-TextCell $battingExtra 19 6 'NO'
Set-TextCell $battingExtra 20 1 '4608'
Set-TextCell $battingExtra 20 2 ''
Set-TextCell $battingExtra 20 3 ''
Set-TextCell $battingExtra 20 4 ''
Set-TextCell $battingExtra 20 5 ''
Set-TextCell $battingExtra 20 6 ''
Set-TextCell $battingExtra 21 1 '4614'
Set-TextCell $battingExtra 21 2 ''
Set-TextCell $battingExtra 21 3 ''
Set-TextCell $battingExtra 21 4 ''
Set-TextCell $battingExtra 21 5 ''
Set-TextCell $battingExtra 21 6 ''

# ---------------------------------------------------------------
# 6. Restore the active sheet/tab to the first sheet
# ---------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
